$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# D-column cells are forced to Text format first so values like "1.001",
# "123.00" or "0.00001118" are preserved verbatim instead of Excel
# auto-coercing them into numbers (which would drop formatting / trailing zeros,
# or switch to scientific notation).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.106.39"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.903.65"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.63"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5044"
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3934"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09580"
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.05"
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.385"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.86"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.890.50"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.332"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001118"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.24"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06604"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.87"
$ws.Range("E20").Value = "  +2.78%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.203"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.153.36"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.306"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.660"
$ws.Range("E26").Value = "  +4.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.109.73"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.83"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.70"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.58"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.084"
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1062"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.619"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.615"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.572"
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06601"
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02425"
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.234"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2184"
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.272"
$ws.Range("E40").Value = "  +8.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.998"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6344"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.35"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.27"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5985"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.725"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.274"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.021"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.00"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.184"
$ws.Range("E51").Value = "  -0.99%  "
